$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.621.22'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '3.514.88'
$ws.Range('E3').Value = '  -1.48%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'586.94"
$ws.Range('E5').Value = '  -2.29%  '
$ws.Range('D6').Value = "'133.00"
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('D7').Value = '3.516.09'
$ws.Range('E7').Value = '  -1.39%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.490"
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('D11').Value = "'7.17"
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').Value = "'0.390"
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '4.111.60'
$ws.Range('E13').Value = '  -1.09%  '
$ws.Range('D14').Value = "'27.94"
$ws.Range('E14').Value = '  +3.32%  '
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('D17').Value = '3.512.98'
$ws.Range('D18').Value = '64.620.08'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').Value = "'5.71"
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('D22').Value = "'393.46"
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').Value = "'0.580"
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '3.655.33'
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('D29').Value = "'7.53"
$ws.Range('E29').Value = '  -4.14%  '
$ws.Range('D30').Value = "'1.00"
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('D32').Value = "'8.25"
$ws.Range('E32').Value = '  -4.34%  '
$ws.Range('D33').Value = '3.518.79'
$ws.Range('E33').Value = '  -1.33%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = "'24.06"
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('D36').Value = "'0.146"
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').Value = "'1.62"
$ws.Range('E37').Value = '  +5.03%  '
$ws.Range('D38').Value = "'5.29"
$ws.Range('E38').Value = '  +4.29%  '
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('E40').Value = '  +0.37%  '
$ws.Range('D41').Value = "'0.0814"
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = "'26.64"
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = "'0.815"
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').Value = "'42.33"
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('D46').Value = "'1.22"
$ws.Range('E46').Value = '  -3.42%  '
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('D48').Value = "'1.66"
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('D49').Value = '2.472.66'
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').Value = "'6.91"
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('D51').Value = "'0.912"
$ws.Range('E51').Value = '  +4.97%  '
